$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: rename "iaest-measure:*" annotations to "iaest-dimension:*" ---
# (sector-descripcion, mes-nombre, sexo and mes-y-ano move from "medida"/measure
#  to "dim"/dimension in this update)
$ws.Range("C3").Value = "iaest-dimension:sector-descripcion"
$ws.Range("D3").Value = "iaest-dimension:mes-nombre"
$ws.Range("H3").Value = "iaest-dimension:sexo"
$ws.Range("M3").Value = "iaest-dimension:mes-y-ano"

# --- Row 4: "medida" -> "dim" for the same four columns ---
$ws.Range("C4").Value = "dim"
$ws.Range("D4").Value = "dim"
$ws.Range("H4").Value = "dim"
$ws.Range("M4").Value = "dim"

# --- Row 5: data-type annotations updated; the three re-typed columns become
#     "skos:Concept" (they now carry a controlled vocabulary / mapping file) ---
$ws.Range("C5").Value = "skos:Concept"
$ws.Range("D5").Value = "skos:Concept"
$ws.Range("H5").Value = "skos:Concept"
$ws.Range("M5").Value = "xsd:string"

# --- Row 6 (new): mapping file references for the re-typed columns ---
$ws.Range("C6").Value = "mapping-sector-descripcion.xlsx"
$ws.Range("D6").Value = "mapping-mes-nombre.xlsx"
$ws.Range("H6").Value = "mapping-sexo.xlsx"

# Match the formatting used by the rest of the data table (row 5) so the
# new row carries the same cell style as every other data cell.
$ws.Range("C5").Copy()
$ws.Range("C6:D6").PasteSpecial(-4122)
$ws.Range("H5").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
